# Refresh cryptos price/volume snapshot (values scraped upstream as plain text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.807.91"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "2.494.50"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "`'485.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.63%  "
$ws.Range("D6").Value = "`'145.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.47%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "`'0.508"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.57%  "
$ws.Range("D9").Value = "2.508.66"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "`'5.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "`'0.0970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "`'0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "2.917.72"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "55.824.52"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").Value = "`'20.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.49%  "
$ws.Range("E17").Value = "  +3.52%  "
$ws.Range("D18").Value = "2.504.96"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").Value = "`'4.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").Value = "`'10.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.52%  "
$ws.Range("D21").Value = "`'319.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "`'5.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.07%  "
$ws.Range("D24").Value = "`'58.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "`'0.409"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.98%  "
$ws.Range("E26").Value = "  +6.74%  "
$ws.Range("D27").Value = "`'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "2.602.95"
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").Value = "`'7.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.65%  "
$ws.Range("D30").Value = "0.0₃0782"
$ws.Range("E30").Value = "  +10.09%  "
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "`'148.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "`'18.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "`'1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("D35").Value = "`'5.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  +8.65%  "
$ws.Range("D37").Value = "`'3.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("D38").Value = "`'0.862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.42%  "
$ws.Range("D39").Value = "`'34.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").Value = "`'3.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "`'0.612"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "`'0.993"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "`'0.0552"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("D44").Value = "`'1.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.64%  "
$ws.Range("D45").Value = "`'4.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.21%  "
$ws.Range("D46").Value = "`'260.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.60%  "
$ws.Range("D47").Value = "`'10.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "`'0.0227"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("D49").Value = "`'0.0903"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("D50").Value = "1.918.43"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("D51").Value = "`'17.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.12%  "
